$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two rows (7 and 8) refer to the same bird observation report split into
# two records that got mismatched: their Id, coordinates (Ost/Nord), dates
# (Startdatum/Slutdatum) and Observatörer values need to be swapped between
# row 7 and row 8.
$xlPasteValues = -4163

$cols = @("A", "Q", "R", "Y", "AA", "AX")

# Scratch cell (outside the used range) used to hold one side of the swap.
# Copy + PasteSpecial(values) preserves each cell's original data type
# (number vs text) without touching cell formatting/style - a straight
# Value/Value2 assignment of a date-like string such as "2022-03-18" would
# otherwise get auto-converted into a date serial number.
$scratch = $ws.Range("ZZ1")

foreach ($col in $cols) {
    $cell7 = $ws.Range($col + "7")
    $cell8 = $ws.Range($col + "8")

    $cell7.Copy()
    $scratch.PasteSpecial($xlPasteValues)

    $cell8.Copy()
    $cell7.PasteSpecial($xlPasteValues)

    $scratch.Copy()
    $cell8.PasteSpecial($xlPasteValues)
}

$scratch.ClearContents()
$excel.CutCopyMode = $false
